$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 70, pushing all existing rows (70..145) down to (71..146)
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new data record
$ws.Cells.Item(70, 1).Value = 4
$ws.Cells.Item(70, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(70, 3).Value = "Los Lagos"
$ws.Cells.Item(70, 4).Value = 44601
$ws.Cells.Item(70, 5).Value = 10
$ws.Cells.Item(70, 6).Value = 100112009
$ws.Cells.Item(70, 7).Value = "Acelga"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 20
$ws.Cells.Item(70, 11).Value = 10000
$ws.Cells.Item(70, 12).Value = 10000
$ws.Cells.Item(70, 13).Value = 10000
$ws.Cells.Item(70, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(70, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(70, 16).Value = 833
$ws.Cells.Item(70, 17).Value = 12
$ws.Cells.Item(70, 18).Value = "Hortaliza"
